# The sheet had a stray header cell in row 1 and a blank row 2-3 before the
# real header row (row 4) and the data (rows 5-24). Remove the first three
# rows so the real header becomes row 1 and the data shifts up to rows 2-21.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("1:3").Delete()

# Reset the view: select the full table range, which also clears the
# previous scroll position / stale selection left over from before the
# rows were removed.
$ws.Range("A1:H22").Select() | Out-Null
